$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Ghost-Lit Raider', ['{2}{R}', 'Creature — Spirit', '{2}{R}, {T}: Ghost-Lit Raider deals 2 damage to target creature.', 'Channel — {3}{R}, Discard Ghost-Lit Raider: It deals 4 damage to target creature.', '2/1'])"
$ws.Range("A3").Value = "('Kiyomaro, First to Stand', ['{3}{W}{W}', 'Legendary Creature — Spirit', 'Kiyomaro, First to Stand’s power and toughness are each equal to the number of cards in your hand.', 'As long as you have four or more cards in hand, Kiyomaro has vigilance.', 'Whenever Kiyomaro deals damage, if you have seven or more cards in hand, you gain 7 life.', '*/*'])"

$ws.Rows("4:14").Delete()
